$p = $ppt.ActivePresentation
$s = $p.Slides.Add(5, 12)
Write-Output "Count=$($p.Slides.Count)"
Write-Output "Layout=$($s.Layout)"
Write-Output "Shapes=$($s.Shapes.Count)"
